# Applies the "LevelUP e Limite de Atributos" edit:
#  - centers the two bold headings ("Level UP" and "Limite de Atributos")
#  - reword the XP paragraph (xp -> experiencia) and split the "EX:" example
#    into its own "Exemplo:" paragraph
#  - moves the hidden _GoBack bookmark from the last paragraph up onto the
#    "1 ponto de atributo interpretativo por nivel" paragraph
#  - merges a couple of runs that had been split for no semantic reason
#
# Implemented from the bottom of the document upwards so that paragraph
# indices encountered later in the script (earlier in the doc) are not
# disturbed by paragraph-count changes made earlier in the script.

$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerWordXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerWordXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- Paragraph 10 (last paragraph): "O sistema de limite ..." ------------
# Merge the two runs into one and drop the _GoBack bookmark (it moves to
# paragraph 6 below). Replacing the range of the very last paragraph in the
# body leaves behind a spare trailing empty paragraph, so clean that up too.
Set-ParagraphXml 10 '<w:p><w:r><w:t>O sistema de limite se encontra na ficha Evolu&#231;&#227;o dos Personagens. O limite de atributo serve o prop&#243;sito de balancear o sistema contra personagens que tenham um foco muito alto em um atributo s&#243;.</w:t></w:r></w:p>'

$n = $d.Paragraphs.Count
if ($n -gt 10) {
    $prev = $d.Paragraphs($n - 1)
    $last = $d.Paragraphs($n)
    $mergeRange = $d.Range($prev.Range.End - 1, $last.Range.End)
    $mergeRange.Delete()
}

# --- Paragraph 9: "Limite de Atributos" heading -> centered --------------
$d.Paragraphs(9).Range.ParagraphFormat.Alignment = 1

# --- Paragraph 7: "+100 pontos ..." -> single run ------------------------
Set-ParagraphXml 7 '<w:p><w:r><w:t>+100 pontos para distribuir entre as barras de P.V, Chakra e Estamina.</w:t></w:r></w:p>'

# --- Paragraph 6: "1 ponto de atributo interpretativo por nivel" ---------
# Add the (hidden) _GoBack bookmark at the very start of the paragraph.
$p6 = $d.Paragraphs(6)
$bmRange = $p6.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Paragraph 3: the long XP paragraph -> reworded + split in two -------
$para3Xml = '<w:p><w:r><w:t xml:space="preserve">Ap&#243;s completar a barra de </w:t></w:r>' +
    '<w:r><w:t>experi&#234;ncia</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, o personagem ir&#225; evoluir 1 n&#237;vel ent&#227;o a barra &#233; resetada e </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">a experi&#234;ncia </w:t></w:r>' +
    '<w:r><w:t>m&#225;xim</w:t></w:r>' +
    '<w:r><w:t>a</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> &#233; aumentad</w:t></w:r>' +
    '<w:r><w:t>a</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> em +100 por n&#237;vel</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:r><w:t>E</w:t></w:r>' +
    '<w:r><w:t>xemplo</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: n&#237;vel 3 = 300 </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>xp</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> total, n&#237;vel 4 </w:t></w:r>' +
    '<w:r><w:t>= 400</w:t></w:r>' +
    '<w:r><w:t>...</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 3 $para3Xml

# --- Paragraph 2: "Level UP" heading -> centered + spell-check artifact --
$para2Xml = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Level</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> UP</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 2 $para2Xml
